# Updated BoM with a cheaper connector from Newark (Molex -> Multicomp terminal
# blocks, sourced from Newark instead of Mouser).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: 9-circuit terminal block header (was Molex/Mouser 39502-1009) ---
$ws.Range("B9").Value = "Terminal Block Header, 9, 150 V, 12 A, 3.5 mm, Through Hole Right Angle, Header"
$ws.Range("D9").Value = "MULTICOMP"
$ws.Range("E9").Value = "MCTE-03A09"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "Newark"
$ws.Range("G9").Value = "54T7648"
$ws.Range("H9").Value = 0.774

# --- Row 10: 9-circuit pluggable terminal block plug (was Molex/Mouser 39500-0009) ---
$ws.Range("D10").Value = "MULTICOMP"
$ws.Range("E10").Value = "MCTC-52A09"
$ws.Range("F10").Value = "Newark"
$ws.Range("G10").Value = "54T7625"
$ws.Range("H10").Value = 1.27

# Move the selection to reflect where the edit was made.
$ws.Range("H11").Select()
